$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add browser details: set C2 to "Completed" (plain, unstyled cell)
$ws.Cells.Item(2, 3).Value = "Completed"
